# Apply crypto price/volume updates to match the new commit snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Staging cell used to push number-like text (e.g. '1.003') into the sheet
# as TEXT via Copy + PasteSpecial(xlPasteValues), so Excel's automatic
# number/date coercion on Range.Value doesn't silently turn strings like
# '330.10' or '1.165' into numeric values.
$stage = $ws.Range('Z1')
$stage.NumberFormat = '@'

$ws.Range('D2').Value = '30.509.32'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '2.106.68'
$ws.Range('E3').Value = '  +4.69%  '
$ws.Range('E4').Value = '  +0.11%  '
$stage.Value = '330.10'
$stage.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('E6').Value = '  +0.13%  '
$stage.Value = '0.5255'
$stage.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('E8').Value = '  +2.96%  '
$stage.Value = '0.08890'
$stage.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +1.63%  '
$stage.Value = '49.52'
$stage.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +13.70%  '
$stage.Value = '1.165'
$stage.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +2.69%  '
$stage.Value = '24.86'
$stage.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').Value = '2.108.76'
$ws.Range('E13').Value = '  +4.54%  '
$stage.Value = '6.747'
$stage.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.55%  '
$stage.Value = '7.764'
$stage.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +4.00%  '
$stage.Value = '96.58'
$stage.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +2.45%  '
$stage.Value = '1.004'
$stage.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('E18').Value = '  +1.59%  '
$stage.Value = '0.06647'
$stage.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +1.49%  '
$stage.Value = '19.32'
$stage.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +2.43%  '
$ws.Range('E21').Value = '  +0.05%  '
$stage.Value = '6.319'
$stage.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('D23').Value = '30.565.98'
$ws.Range('E23').Value = '  +0.40%  '
$stage.Value = '12.28'
$stage.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +3.75%  '
$stage.Value = '2.344'
$stage.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +4.08%  '
$ws.Range('D26').Value = '2.348.79'
$ws.Range('E26').Value = '  +4.22%  '
$stage.Value = '22.48'
$stage.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +0.02%  '
$stage.Value = '2.635'
$stage.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +7.13%  '
$stage.Value = '162.32'
$stage.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +0.16%  '
$stage.Value = '133.00'
$stage.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +1.21%  '
$stage.Value = '1.224'
$stage.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +7.46%  '
$stage.Value = '0.1073'
$stage.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +1.83%  '
$stage.Value = '1.688'
$stage.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +23.75%  '
$stage.Value = '6.240'
$stage.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +2.53%  '
$stage.Value = '3.897'
$stage.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +1.76%  '
$stage.Value = '10.21'
$stage.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +10.79%  '
$stage.Value = '0.02590'
$stage.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$stage.Value = '0.06752'
$stage.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$stage.Value = '5.518'
$stage.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +1.05%  '
$stage.Value = '12.74'
$stage.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +2.37%  '
$stage.Value = '0.2287'
$stage.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +3.18%  '
$stage.Value = '0.6932'
$stage.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +3.97%  '
$stage.Value = '1.277'
$stage.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +3.01%  '
$stage.Value = '1.002'
$stage.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +0.07%  '
$stage.Value = '0.6435'
$stage.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +4.09%  '
$stage.Value = '14.08'
$stage.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +2.80%  '
$stage.Value = '2.230'
$stage.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +1.16%  '
$stage.Value = '3.638'
$stage.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  -0.56%  '
$stage.Value = '1.219'
$stage.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +10.02%  '
$stage.Value = '82.88'
$stage.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +2.00%  '

# Clean up the staging cell/clipboard so no stray content or used-range
# expansion is left behind.
$stage.Clear()
$excel.CutCopyMode = 0

